$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 2")
$ws.Activate()
$ws.Range("A16").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("F41").Select()
